# Edit script: apply weekly price update for Chirimoya (Terminal La Palmera de La Serena)
# - Rows 28-30: new week of data (fecha 44466), Provincia del Elqui figures shift in
# - Rows 31-41: existing rows' data cascades (dates/qualities/prices rotate down)
# - Rows 42-43: updated in place
# - Rows 44-46: new rows appended (previously "dropped" oldest week's data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed (unchanged) descriptive columns shared by every data row in this sheet.
$fixedA = 8
$fixedB = "Terminal La Palmera de La Serena"
$fixedC = "Coquimbo"
$fixedE = 4
$fixedF = "Fruta"
$fixedG = 100107
$fixedH = "Otros"
$fixedI = 100107002
$fixedJ = "Chirimoya"
$fixedK = "Cultivar IV Región"

# Per-row target values: Row, D(Fecha), L(Calidad), M(Volumen), N(Precio min), O(Precio max),
#                         P(Precio prom pond), Q(Unidad comerc.), R(Origen), S(Precio $/Kg), T(Kg/unidad)
$rows = @(
    @(28, 44466, 'Especial', 160, 2100, 2200, 2150, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2150, 1),
    @(29, 44466, 'Primera', 240, 1700, 1800, 1750, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 1750, 1),
    @(30, 44466, 'Segunda', 200, 1300, 1400, 1350, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 1350, 1),
    @(31, 44463, 'Especial', 240, 2600, 2700, 2650, '$/kilo (en caja de 15 kilos)', 'Provincia de Limarí', 2650, 1),
    @(32, 44463, 'Primera', 300, 2200, 2300, 2250, '$/kilo (en caja de 15 kilos)', 'Provincia de Limarí', 2250, 1),
    @(33, 44463, 'Segunda', 240, 1900, 2000, 1950, '$/kilo (en caja de 15 kilos)', 'Provincia de Limarí', 1950, 1),
    @(34, 44166, 'Especial', 300, 14000, 14500, 14250, '$/bandeja 8 kilos', 'Provincia de Limarí', 1781, 8),
    @(35, 44166, 'Primera', 200, 12000, 12500, 12250, '$/bandeja 8 kilos', 'Provincia de Limarí', 1531, 8),
    @(36, 44168, 'Especial', 240, 14000, 14500, 14250, '$/bandeja 8 kilos', 'Provincia de Limarí', 1781, 8),
    @(37, 44168, 'Primera', 200, 12000, 12500, 12250, '$/bandeja 8 kilos', 'Provincia de Limarí', 1531, 8),
    @(38, 44168, 'Segunda', 200, 9500, 10000, 9750, '$/bandeja 8 kilos', 'Provincia de Limarí', 1219, 8),
    @(39, 44162, 'Especial', 340, 14000, 14500, 14250, '$/bandeja 8 kilos', 'Provincia de Limarí', 1781, 8),
    @(40, 44162, 'Primera', 300, 12000, 12500, 12250, '$/bandeja 8 kilos', 'Provincia de Limarí', 1531, 8),
    @(41, 44162, 'Segunda', 200, 9500, 10000, 9750, '$/bandeja 8 kilos', 'Provincia de Limarí', 1219, 8),
    @(42, 44410, 'Primera', 240, 2400, 2500, 2450, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2450, 1),
    @(43, 44410, 'Segunda', 240, 2000, 2100, 2050, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2050, 1),
    @(44, 44410, 'Tercera', 200, 1600, 1700, 1650, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 1650, 1),
    @(45, 44411, 'Primera', 600, 2400, 2500, 2450, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2450, 1),
    @(46, 44411, 'Segunda', 400, 2000, 2100, 2050, '$/kilo (en caja de 15 kilos)', 'Provincia del Elquí', 2050, 1)
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    # Ensure descriptive columns are present/correct for every row (new rows 44-46 start empty).
    $ws.Cells.Item($rowNum, 1).Value = $fixedA
    $ws.Cells.Item($rowNum, 2).Value = $fixedB
    $ws.Cells.Item($rowNum, 3).Value = $fixedC
    $ws.Cells.Item($rowNum, 4).Value = $r[1]
    $ws.Cells.Item($rowNum, 5).Value = $fixedE
    $ws.Cells.Item($rowNum, 6).Value = $fixedF
    $ws.Cells.Item($rowNum, 7).Value = $fixedG
    $ws.Cells.Item($rowNum, 8).Value = $fixedH
    $ws.Cells.Item($rowNum, 9).Value = $fixedI
    $ws.Cells.Item($rowNum, 10).Value = $fixedJ
    $ws.Cells.Item($rowNum, 11).Value = $fixedK
    $ws.Cells.Item($rowNum, 12).Value = $r[2]
    $ws.Cells.Item($rowNum, 13).Value = $r[3]
    $ws.Cells.Item($rowNum, 14).Value = $r[4]
    $ws.Cells.Item($rowNum, 15).Value = $r[5]
    $ws.Cells.Item($rowNum, 16).Value = $r[6]
    $ws.Cells.Item($rowNum, 17).Value = $r[7]
    $ws.Cells.Item($rowNum, 18).Value = $r[8]
    $ws.Cells.Item($rowNum, 19).Value = $r[9]
    $ws.Cells.Item($rowNum, 20).Value = $r[10]

    # Column D (Fecha) keeps the date/time custom number format used throughout the sheet.
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
